$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 221.5
$ws.Range("I2").Value = 156.7
$ws.Range("J2").Value = 545.5
$ws.Range("K2").Value = 156.7
$ws.Range("L2").Value = 545.5
$ws.Range("M2").Value = -43.69999999999999
$ws.Range("N2").Value = -771.5

$ws.Range("H32").Value = 1699
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1699
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1699
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2351

$ws.Range("H40").Value = 4422.5386
$ws.Range("I40").Value = 3132.5
$ws.Range("J40").Value = 5528.2856
$ws.Range("K40").Value = 3132.5
$ws.Range("L40").Value = 5528.2856
$ws.Range("M40").Value = -2957.5
$ws.Range("N40").Value = -5878.2856

$ws.Range("H42").Value = 756.1429000000001
$ws.Range("I42").Value = 17
$ws.Range("J42").Value = 1051.8
$ws.Range("K42").Value = 51
$ws.Range("L42").Value = 3155.4
$ws.Range("M42").Value = 179
$ws.Range("N42").Value = -3615.4

$ws.Range("H53").Value = 354.2857
$ws.Range("I53").Value = 334.5
$ws.Range("J53").Value = 369.125
$ws.Range("K53").Value = 334.5
$ws.Range("L53").Value = 369.125
$ws.Range("M53").Value = 302.5
$ws.Range("N53").Value = -1643.125

$ws.Range("H64").Value = 14705.294
$ws.Range("I64").Value = 8749.375
$ws.Range("J64").Value = 19999.445
$ws.Range("K64").Value = 8749.375
$ws.Range("L64").Value = 19999.445
$ws.Range("M64").Value = -8501.375
$ws.Range("N64").Value = -20495.445

$ws.Range("H67").Value = 14705.294
$ws.Range("I67").Value = 8749.375
$ws.Range("J67").Value = 19999.445
$ws.Range("K67").Value = 8749.375
$ws.Range("L67").Value = 19999.445
$ws.Range("M67").Value = -7891.375
$ws.Range("N67").Value = -21715.445

$ws.Range("H100").Value = 2005.3572
$ws.Range("I100").Value = 1914
$ws.Range("K100").Value = 1914
$ws.Range("M100").Value = -1373

$ws.Range("H118").Value = 966.3333
$ws.Range("I118").Value = 966.3333
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2898.9999
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -1241.9999
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1408.625
$ws.Range("I2").Value = 179.42857
$ws.Range("K2").Value = 179.42857
$ws.Range("M2").Value = -66.42857000000001

$ws.Range("H45").Value = 1637.25
$ws.Range("I45").Value = 1562.2
$ws.Range("J45").Value = 2012.5
$ws.Range("K45").Value = 1562.2
$ws.Range("L45").Value = 2012.5
$ws.Range("M45").Value = -1185.2
$ws.Range("N45").Value = -2766.5

$ws.Range("H74").Value = 2496
$ws.Range("I74").Value = 2496
$ws.Range("K74").Value = 2496
$ws.Range("M74").Value = -1622

$ws.Range("H77").Value = 2496
$ws.Range("I77").Value = 2496
$ws.Range("K77").Value = 12480
$ws.Range("M77").Value = -8112

$ws.Range("H116").Value = 1408.625
$ws.Range("I116").Value = 179.42857
$ws.Range("K116").Value = 179.42857
$ws.Range("M116").Value = 2114.57143

$ws.Range("H119").Value = 74000
$ws.Range("J119").Value = 74000
$ws.Range("L119").Value = 74000
$ws.Range("N119").Value = -83676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1408.625
$ws.Range("I3").Value = 179.42857
$ws.Range("K3").Value = 179.42857
$ws.Range("M3").Value = -65.42857000000001

$ws.Range("H87").Value = 53777.5
$ws.Range("I87").Value = 60000
$ws.Range("J87").Value = 47555
$ws.Range("K87").Value = 60000
$ws.Range("L87").Value = 47555
$ws.Range("M87").Value = -58752
$ws.Range("N87").Value = -50051

$ws.Range("H90").Value = 53777.5
$ws.Range("I90").Value = 60000
$ws.Range("J90").Value = 47555
$ws.Range("K90").Value = 180000
$ws.Range("L90").Value = 142665
$ws.Range("M90").Value = -173760
$ws.Range("N90").Value = -155145

$ws.Range("H96").Value = 13499
$ws.Range("I96").Value = 13499
$ws.Range("K96").Value = 13499
$ws.Range("M96").Value = -10753

$ws.Range("H105").Value = 1611
$ws.Range("I105").Value = 1611
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1611
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 136
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 12800
$ws.Range("I60").Value = 8500
$ws.Range("J60").Value = 30000
$ws.Range("K60").Value = 8500
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -7989
$ws.Range("N60").Value = -31022

$ws.Range("H134").Value = 1781.8948
$ws.Range("I134").Value = 1286.0588
$ws.Range("K134").Value = 3858.1764
$ws.Range("M134").Value = -1323.1764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 91.666664
$ws.Range("I26").Value = 91.666664
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 274.999992
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 13.00000799999998
$ws.Range("N26").ClearContents()

$ws.Range("H40").Value = 164.2
$ws.Range("I40").Value = 189.75
$ws.Range("J40").Value = 62
$ws.Range("K40").Value = 759
$ws.Range("L40").Value = 248
$ws.Range("M40").Value = -690
$ws.Range("N40").Value = -386

$ws.Range("H64").Value = 1250
$ws.Range("I64").Value = 1250
$ws.Range("K64").Value = 3750
$ws.Range("M64").Value = -3480

$ws.Range("H67").Value = 1250
$ws.Range("I67").Value = 1250
$ws.Range("K67").Value = 3750
$ws.Range("M67").Value = -2814

$ws.Range("H80").Value = 2170.4
$ws.Range("I80").Value = 2250.5
$ws.Range("J80").Value = 1850
$ws.Range("K80").Value = 6751.5
$ws.Range("L80").Value = 5550
$ws.Range("M80").Value = -5815.5
$ws.Range("N80").Value = -7422

$ws.Range("H83").Value = 2170.4
$ws.Range("I83").Value = 2250.5
$ws.Range("J83").Value = 1850
$ws.Range("K83").Value = 20254.5
$ws.Range("L83").Value = 16650
$ws.Range("M83").Value = -15574.5
$ws.Range("N83").Value = -26010

$ws.Range("H122").Value = 955.8
$ws.Range("I122").Value = 632
$ws.Range("J122").Value = 1171.6666
$ws.Range("K122").Value = 5688
$ws.Range("L122").Value = 10544.9994
$ws.Range("M122").Value = -3238
$ws.Range("N122").Value = -15444.9994

$ws.Range("H130").Value = 10000
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 500005000
$ws.Range("J126").Value = 9999
$ws.Range("L126").Value = 29997
$ws.Range("N126").Value = -34937

$ws.Range("H132").Value = 3438.7
$ws.Range("I132").Value = 2956.5
$ws.Range("K132").Value = 8869.5
$ws.Range("M132").Value = -6339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251

$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = -21256

$ws.Range("H136").Value = 4686.875
$ws.Range("I136").Value = 4499.1665
$ws.Range("J136").Value = 5250
$ws.Range("K136").Value = 13497.4995
$ws.Range("L136").Value = 15750
$ws.Range("M136").Value = -10947.4995
$ws.Range("N136").Value = -20850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4499.5
$ws.Range("I136").Value = 4333.3335
$ws.Range("K136").Value = 13000.0005
$ws.Range("M136").Value = -10450.0005
